# Insert a new weekly observation row for "Ciboulette" (Vega Central
# Mapocho de Santiago) above the existing row 521. Excel's Rows(...).Insert()
# shifts row 521 and everything below it down by one (521 -> 522, ...,
# 558 -> 559), matching the rest of the dataset exactly, and we only need
# to populate the freshly inserted row with its new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(521).Insert()

$ws.Range("A521").Value = 9
$ws.Range("B521").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C521").Value = "Metropolitana"
$ws.Range("D521").Value = 45021
$ws.Range("E521").Value = 13
$ws.Range("F521").Value = 100112039
$ws.Range("G521").Value = "Ciboulette"
$ws.Range("H521").Value = "Sin especificar"
$ws.Range("I521").Value = "Primera"
$ws.Range("J521").Value = 430
$ws.Range("K521").Value = 1000
$ws.Range("L521").Value = 1200
$ws.Range("M521").Value = 1116
$ws.Range("N521").Value = "`$/docena de atados"
$ws.Range("O521").Value = "Región Metropolitana"
$ws.Range("P521").Value = 372
$ws.Range("Q521").Value = 3
$ws.Range("R521").Value = "Hortaliza"
